$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.696.65"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.469.22"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.480"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("D12").Value = "4.064.41"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.35%  "
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "3.459.21"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "62.814.04"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +5.28%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -9.20%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("E37").Value = "  +7.34%  "
$ws.Range("E38").Value = "  +21.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "168.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "3.510.15"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("D47").Value = "2.605.91"
$ws.Range("E47").Value = "  +6.28%  "
$ws.Range("E48").Value = "  +11.25%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
